$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2733
$ws.Range("I40").Value = 2733
$ws.Range("K40").Value = 2733
$ws.Range("M40").Value = -2558

$ws.Range("H68").Value = 463678.34
$ws.Range("J68").Value = 195517.5
$ws.Range("L68").Value = 195517.5
$ws.Range("N68").Value = -197015.5

$ws.Range("H71").Value = 463678.34
$ws.Range("J71").Value = 195517.5
$ws.Range("L71").Value = 586552.5
$ws.Range("N71").Value = -594040.5

$ws.Range("H97").Value = 1963.5
$ws.Range("I97").Value = 539
$ws.Range("K97").Value = 1617
$ws.Range("M97").Value = -1121

$ws.Range("H100").Value = 2714.6667
$ws.Range("I100").Value = 2717.6
$ws.Range("J100").Value = 2700
$ws.Range("K100").Value = 2717.6
$ws.Range("L100").Value = 2700
$ws.Range("M100").Value = -2176.6
$ws.Range("N100").Value = -3782

$ws.Range("H113").Value = 8220.1
$ws.Range("I113").Value = 7466.8887
$ws.Range("K113").Value = 7466.8887
$ws.Range("M113").Value = -4212.8887

$ws.Range("H116").Value = 5950
$ws.Range("I116").Value = 5900
$ws.Range("J116").Value = 6000
$ws.Range("K116").Value = 5900
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = -2458
$ws.Range("N116").Value = -12884

$ws.Range("H138").Value = 1991.4255
$ws.Range("J138").Value = 2246.353
$ws.Range("L138").Value = 6739.059
$ws.Range("N138").Value = -17019.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2900.7407
$ws.Range("I61").Value = 1875.75
$ws.Range("K61").Value = 1875.75
$ws.Range("M61").Value = -1663.75

$ws.Range("H88").Value = 1929.8667
$ws.Range("I88").Value = 1468
$ws.Range("K88").Value = 1468
$ws.Range("M88").Value = -1062

$ws.Range("H91").Value = 1929.8667
$ws.Range("I91").Value = 1468
$ws.Range("K91").Value = 1468
$ws.Range("M91").Value = -64

$ws.Range("H110").Value = 1814.4231
$ws.Range("I110").Value = 1618.8096
$ws.Range("K110").Value = 1618.8096
$ws.Range("M110").Value = 426.1904

$ws.Range("H136").Value = 2900.7407
$ws.Range("I136").Value = 1875.75
$ws.Range("K136").Value = 5627.25
$ws.Range("M136").Value = -3077.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1556.9445
$ws.Range("I20").Value = 1369.0667
$ws.Range("K20").Value = 1369.0667
$ws.Range("M20").Value = -1122.0667

$ws.Range("H86").Value = 3670.65
$ws.Range("I86").Value = 3678
$ws.Range("K86").Value = 3678
$ws.Range("M86").Value = -2555

$ws.Range("H89").Value = 3670.65
$ws.Range("I89").Value = 3678
$ws.Range("K89").Value = 18390
$ws.Range("M89").Value = -12774

$ws.Range("H94").Value = 824.84
$ws.Range("I94").Value = 424.77777
$ws.Range("J94").Value = 1853.5714
$ws.Range("K94").Value = 424.77777
$ws.Range("L94").Value = 1853.5714
$ws.Range("M94").Value = 26.22223000000002
$ws.Range("N94").Value = -2755.5714

$ws.Range("H107").Value = 2736.1538
$ws.Range("I107").Value = 2690.08
$ws.Range("K107").Value = 2690.08
$ws.Range("M107").Value = -770.0799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2324.05
$ws.Range("I31").Value = 2037.1111
$ws.Range("K31").Value = 2037.1111
$ws.Range("M31").Value = -1742.1111

$ws.Range("H34").Value = 2324.05
$ws.Range("I34").Value = 2037.1111
$ws.Range("K34").Value = 2037.1111
$ws.Range("M34").Value = -1835.1111

$ws.Range("H35").Value = 295
$ws.Range("I35").Value = 295
$ws.Range("K35").Value = 295
$ws.Range("M35").Value = -1

$ws.Range("H58").Value = 6046.4375
$ws.Range("I58").Value = 2949.7144
$ws.Range("J58").Value = 8455
$ws.Range("K58").Value = 2949.7144
$ws.Range("L58").Value = 8455
$ws.Range("M58").Value = -2746.7144
$ws.Range("N58").Value = -8861

$ws.Range("H122").Value = 2076.8
$ws.Range("I122").Value = 1696.5
$ws.Range("J122").Value = 2330.3333
$ws.Range("K122").Value = 5089.5
$ws.Range("L122").Value = 6990.999899999999
$ws.Range("M122").Value = -2639.5
$ws.Range("N122").Value = -11890.9999

$ws.Range("H136").Value = 6046.4375
$ws.Range("I136").Value = 2949.7144
$ws.Range("J136").Value = 8455
$ws.Range("K136").Value = 8849.143199999999
$ws.Range("L136").Value = 25365
$ws.Range("M136").Value = -6299.143199999999
$ws.Range("N136").Value = -30465

$ws.Range("H141").Value = 287500
$ws.Range("J141").Value = 287500
$ws.Range("L141").Value = 287500
$ws.Range("N141").Value = -297860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 865.6667
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 865.6667
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2597.0001
$ws.Range("N22").Value = -2935.0001

$ws.Range("H27").Value = 865.6667
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 865.6667
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2597.0001
$ws.Range("N27").Value = -2801.0001

$ws.Range("H125").Value = 2962
$ws.Range("I125").Value = 2962
$ws.Range("K125").Value = 8886
$ws.Range("M125").Value = -3966

$ws.Range("H129").Value = 9098116
$ws.Range("J129").Value = 16674100
$ws.Range("L129").Value = 50022300
$ws.Range("N129").Value = -50032300

$ws.Range("H130").Value = 3165.5
$ws.Range("I130").Value = 3165.5
$ws.Range("K130").Value = 9496.5
$ws.Range("M130").Value = -4476.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26837.4
$ws.Range("I70").Value = 37977.562
$ws.Range("K70").Value = 37977.562
$ws.Range("M70").Value = -37707.562

$ws.Range("H73").Value = 26837.4
$ws.Range("I73").Value = 37977.562
$ws.Range("K73").Value = 37977.562
$ws.Range("M73").Value = -37041.562

$ws.Range("H102").Value = 3467.8572
$ws.Range("I102").Value = 3554.7
$ws.Range("K102").Value = 3554.7
$ws.Range("M102").Value = -1932.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5290.2
$ws.Range("I40").Value = 5042.3687
$ws.Range("K40").Value = 5042.3687
$ws.Range("M40").Value = -4906.3687

$ws.Range("H46").Value = 1434
$ws.Range("I46").Value = 1276.5
$ws.Range("K46").Value = 1276.5
$ws.Range("M46").Value = -1088.5

$ws.Range("H61").Value = 4414
$ws.Range("I61").Value = 4877.8
$ws.Range("K61").Value = 4877.8
$ws.Range("M61").Value = -4675.8

$ws.Range("H113").Value = 4414
$ws.Range("I113").Value = 4877.8
$ws.Range("K113").Value = 4877.8
$ws.Range("M113").Value = -2707.8

$ws.Range("H132").Value = 5221.7856
$ws.Range("I132").Value = 3201.125
$ws.Range("J132").Value = 7916
$ws.Range("K132").Value = 9603.375
$ws.Range("L132").Value = 23748
$ws.Range("M132").Value = -7073.375
$ws.Range("N132").Value = -28808

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0

$ws.Range("H100").Value = 830
$ws.Range("I100").Value = 830
$ws.Range("K100").Value = 1660
$ws.Range("M100").Value = -1119

$ws.Range("H107").Value = 970.8571
$ws.Range("I107").Value = 960.9231
$ws.Range("K107").Value = 2882.7693
$ws.Range("M107").Value = -962.7692999999999

$wb.Worksheets.Item("CUL").Range("M22").ClearContents()
$wb.Worksheets.Item("CUL").Range("M27").ClearContents()
$wb.Worksheets.Item("WVR").Range("M54").ClearContents()

Write-Host "Applied all changes"